$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the 6bcce489 (second) row
# refreshed as part of regenerating the handback report.
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-05 15:01:16"

# "zh-cn" sheet: Correspond Handoff/Handback datetimes refreshed for the
# 04f690e6 (first) row.
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-05 15:01:01"
$wsZhCn.Range("K2").Value = "2016-09-05 15:01:34"

# "de-de" sheet: Correspond Handback datetime refreshed for the 04f690e6
# (first) row.
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-05 15:01:42"
